# Updated symbol list on Mon Dec 12 05:48:26 UTC 2022 with GitHub Actions
#
# Applies the cell-value changes described by the upstream OOXML diff to
# the "cryptos" worksheet. All data cells on this sheet are stored as
# text (inlineStr) even when their contents look numeric (prices like
# "281.15"), so every write below forces the destination cell to Text
# format before assigning the value and then restores the "Normal" cell
# style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Address, $Value) {
    $cell = $Sheet.Range($Address)
    # Force text storage so numeric-looking strings (e.g. "281.54") are
    # not silently coerced into real numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    # Drop back to the default style so we don't leave a "Text" number
    # format applied to the cell (matches the original, unstyled cells).
    $cell.Style = "Normal"
}

# Row 2 - BNB
Set-TextCell $ws "D2" "281.54"

# Row 3 - OKB
Set-TextCell $ws "D3" "20.70"

# Row 4 - HuobiToken
Set-TextCell $ws "D4" "6.224"

# Row 6 - GateToken
Set-TextCell $ws "D6" "3.582"

# Row 7 - was FTXToken, now KuCoinToken (rows 7 & 8 swap places)
Set-TextCell $ws "B7" "KuCoinToken"
Set-TextCell $ws "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell $ws "D7" "6.567"
Set-TextCell $ws "E7" "6KuCoinTokenKCS"

# Row 8 - was KuCoinToken, now FTXToken
Set-TextCell $ws "B8" "FTXToken"
Set-TextCell $ws "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws "D8" "1.505"
Set-TextCell $ws "E8" "7FTXTokenFTT"

# Row 9 - MXToken
Set-TextCell $ws "D9" "0.8191"

# Row 10 - One
Set-TextCell $ws "D10" "0.01386"

# Row 11 - WazirX
Set-TextCell $ws "D11" "0.1642"

# Row 12 - MandalaExchangeToken
Set-TextCell $ws "D12" "0.08398"

# Row 13 - LiechtensteinCryptoassetsExchange
Set-TextCell $ws "D13" "0.03535"

# Row 15
Set-TextCell $ws "D15" "0.09142"

# Row 16
Set-TextCell $ws "D16" "3.713"

# Row 17
Set-TextCell $ws "D17" "0.001643"

# Row 18
Set-TextCell $ws "D18" "0.04729"

# Row 19
Set-TextCell $ws "D19" "0.006505"

# Row 20
Set-TextCell $ws "D20" "0.006167"

# Row 22
Set-TextCell $ws "D22" "0.0001601"

# Row 23
Set-TextCell $ws "D23" "3.776"

# Row 26
Set-TextCell $ws "D26" "0.1251"

# Row 40 - IDEX
Set-TextCell $ws "D40" "0.04701"

# Row 41 - KickToken
Set-TextCell $ws "D41" "0.007195"

# Row 42 - was BKEXToken, now CEJI (rows 42 & 43 swap places)
Set-TextCell $ws "B42" "CEJI"
Set-TextCell $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell $ws "D42" "0.004503"
Set-TextCell $ws "E42" "41CEJICEJI"

# Row 43 - was CEJI, now BKEXToken
Set-TextCell $ws "B43" "BKEXToken"
Set-TextCell $ws "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws "D43" "0.1098"
Set-TextCell $ws "E43" "42BKEXTokenBKK"

# Row 44 - LocalTraders
Set-TextCell $ws "D44" "0.01108"

# Row 45 - CoinLion
Set-TextCell $ws "D45" "0.00006508"

# Row 48 - BOLO
Set-TextCell $ws "D48" "0.002724"
